# Update the "想去人数" (column F) counts on both the "展览" sheet and the
# "全部类型" sheet, which carry duplicate data for this event list.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F
$updates = @{
    2 = 2283
    3 = 1738
    5 = 1095
    6 = 863
    8 = 5851
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
